$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.291.06"
$ws.Range("E2").Value = "  +2.98%  "

$ws.Range("D3").Value = "3.246.23"
$ws.Range("E3").Value = "  +5.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'576.00"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("D6").Value = "'154.69"
$ws.Range("E6").Value = "  +7.30%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.237.31"
$ws.Range("E8").Value = "  +5.54%  "

$ws.Range("D9").Value = "'0.513"
$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = "  +9.23%  "

$ws.Range("D11").Value = "'0.165"
$ws.Range("E11").Value = "  +4.25%  "

$ws.Range("E12").Value = "  +3.80%  "

$ws.Range("D13").Value = "'37.76"
$ws.Range("E13").Value = "  +4.11%  "

$ws.Range("D14").Value = "'0.0000236"
$ws.Range("E14").Value = "  +4.32%  "

$ws.Range("D15").Value = "3.761.44"
$ws.Range("E15").Value = "  +5.25%  "

$ws.Range("D16").Value = "'563.52"
$ws.Range("E16").Value = "  +13.58%  "

$ws.Range("D17").Value = "66.405.85"
$ws.Range("E17").Value = "  +3.00%  "

$ws.Range("D18").Value = "3.245.59"
$ws.Range("E18").Value = "  +5.14%  "

$ws.Range("E19").Value = "  +2.94%  "

$ws.Range("E20").Value = "  +5.37%  "

$ws.Range("D21").Value = "'14.41"
$ws.Range("E21").Value = "  +3.88%  "

$ws.Range("D22").Value = "'0.745"
$ws.Range("E22").Value = "  +6.85%  "

$ws.Range("D23").Value = "'7.88"
$ws.Range("E23").Value = "  +8.52%  "

$ws.Range("D24").Value = "'13.54"
$ws.Range("E24").Value = "  +4.63%  "

$ws.Range("D25").Value = "'82.21"
$ws.Range("E25").Value = "  +3.68%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = "  +16.91%  "

$ws.Range("E28").Value = "  +5.42%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +6.45%  "

$ws.Range("D30").Value = "'27.99"
$ws.Range("E30").Value = "  +5.75%  "

$ws.Range("E31").Value = "  +2.20%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = "  +4.04%  "

$ws.Range("D34").Value = "'566.93"
$ws.Range("E34").Value = "  +7.68%  "

$ws.Range("D35").Value = "'5.75"
$ws.Range("E35").Value = "  +4.09%  "

$ws.Range("D36").Value = "'6.41"
$ws.Range("E36").Value = "  +5.75%  "

$ws.Range("D37").Value = "'55.33"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").Value = "'0.0454"
$ws.Range("E38").Value = "  +11.66%  "

$ws.Range("D39").Value = "'0.0868"
$ws.Range("E39").Value = "  +7.13%  "

$ws.Range("D40").Value = "'3.07"
$ws.Range("E40").Value = "  +14.46%  "

$ws.Range("E41").Value = "  +5.43%  "

$ws.Range("D42").Value = "3.149.07"
$ws.Range("E42").Value = "  +6.57%  "

$ws.Range("D43").Value = "'8.62"
$ws.Range("E43").Value = "  +1.72%  "

$ws.Range("D44").Value = "'0.275"
$ws.Range("E44").Value = "  +9.79%  "

$ws.Range("D45").Value = "'2.32"
$ws.Range("E45").Value = "  +6.76%  "

$ws.Range("D46").Value = "'26.67"
$ws.Range("E46").Value = "  +3.78%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0560"
$ws.Range("E48").Value = "  +2.84%  "

$ws.Range("E49").Value = "  +3.08%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'123.92"
$ws.Range("E50").Value = "  +3.02%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "  +8.42%  "
